$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3545.4546
$ws.Range("I64").Value = 3542.8572
$ws.Range("J64").Value = 3550
$ws.Range("K64").Value = 3542.8572
$ws.Range("L64").Value = 3550
$ws.Range("M64").Value = -3294.8572
$ws.Range("N64").Value = -4046

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3545.4546
$ws.Range("I67").Value = 3542.8572
$ws.Range("J67").Value = 3550
$ws.Range("K67").Value = 3542.8572
$ws.Range("L67").Value = 3550
$ws.Range("M67").Value = -2684.8572
$ws.Range("N67").Value = -5266

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5251.6665
$ws.Range("I74").Value = 3755
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 3755
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -2819
$ws.Range("N74").Value = -7872

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5251.6665
$ws.Range("I77").Value = 3755
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 18775
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -14095
$ws.Range("N77").Value = -39360

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 977.3692
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 986.3905999999999
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 2959.1718
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -5175.1718

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1438.6
$ws.Range("I137").Value = 1108.9259
$ws.Range("J137").Value = 2123.3076
$ws.Range("K137").Value = 3326.7777
$ws.Range("L137").Value = 6369.9228
$ws.Range("M137").Value = -776.7776999999996
$ws.Range("N137").Value = -11469.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1129.6
$ws.Range("I74").Value = 1006.4186
$ws.Range("J74").Value = 1441.1765
$ws.Range("K74").Value = 1006.4186
$ws.Range("L74").Value = 1441.1765
$ws.Range("M74").Value = -132.4186
$ws.Range("N74").Value = -3189.1765

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1129.6
$ws.Range("I77").Value = 1006.4186
$ws.Range("J77").Value = 1441.1765
$ws.Range("K77").Value = 5032.093
$ws.Range("L77").Value = 7205.8825
$ws.Range("M77").Value = -664.0929999999998
$ws.Range("N77").Value = -15941.8825

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1051579.2
$ws.Range("I132").Value = 954.3421
$ws.Range("J132").Value = 3269565
$ws.Range("K132").Value = 2863.0263
$ws.Range("L132").Value = 9808695
$ws.Range("M132").Value = -333.0263
$ws.Range("N132").Value = -9813755

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 23716
$ws.Range("J55").Value = 23716
$ws.Range("L55").Value = 23716
$ws.Range("N55").Value = -24262

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 16266
$ws.Range("J109").Value = 16266
$ws.Range("L109").Value = 16266
$ws.Range("N109").Value = -19040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3537.4043
$ws.Range("I31").Value = 6142.263
$ws.Range("J31").Value = 1769.8214
$ws.Range("K31").Value = 6142.263
$ws.Range("L31").Value = 1769.8214
$ws.Range("M31").Value = -5847.263
$ws.Range("N31").Value = -2359.8214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3537.4043
$ws.Range("I34").Value = 6142.263
$ws.Range("J34").Value = 1769.8214
$ws.Range("K34").Value = 6142.263
$ws.Range("L34").Value = 1769.8214
$ws.Range("M34").Value = -5940.263
$ws.Range("N34").Value = -2173.8214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 10753875
$ws.Range("I132").Value = 995.84
$ws.Range("J132").Value = 55557540
$ws.Range("K132").Value = 2987.52
$ws.Range("L132").Value = 166672620
$ws.Range("M132").Value = -457.52
$ws.Range("N132").Value = -166677680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 34900
$ws.Range("J137").Value = 34900
$ws.Range("L137").Value = 34900
$ws.Range("N137").Value = -45100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 39058.375
$ws.Range("J141").Value = 39058.375
$ws.Range("L141").Value = 39058.375
$ws.Range("N141").Value = -49418.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2364.25
$ws.Range("I81").Value = 99
$ws.Range("J81").Value = 2687.8572
$ws.Range("K81").Value = 297
$ws.Range("L81").Value = 8063.571599999999
$ws.Range("M81").Value = 826
$ws.Range("N81").Value = -10309.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2364.25
$ws.Range("I84").Value = 99
$ws.Range("J84").Value = 2687.8572
$ws.Range("K84").Value = 891
$ws.Range("L84").Value = 24190.7148
$ws.Range("M84").Value = 4725
$ws.Range("N84").Value = -35422.7148

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1422.1111
$ws.Range("I103").Value = 724.75
$ws.Range("J103").Value = 1980
$ws.Range("K103").Value = 2174.25
$ws.Range("L103").Value = 5940
$ws.Range("M103").Value = -1295.25
$ws.Range("N103").Value = -7698

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 916.25
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 929.1237
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2787.3711
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12867.3711

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3783.0425
$ws.Range("I132").Value = 1203.8055
$ws.Range("J132").Value = 12224.182
$ws.Range("K132").Value = 3611.4165
$ws.Range("L132").Value = 36672.546
$ws.Range("M132").Value = -1081.4165
$ws.Range("N132").Value = -41732.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 21722.777
$ws.Range("I109").Value = 20000
$ws.Range("J109").Value = 21938.125
$ws.Range("K109").Value = 20000
$ws.Range("L109").Value = 21938.125
$ws.Range("M109").Value = -18613
$ws.Range("N109").Value = -24712.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11915.695
$ws.Range("I132").Value = 2761.2856
$ws.Range("J132").Value = 26155.889
$ws.Range("K132").Value = 8283.856800000001
$ws.Range("L132").Value = 78467.667
$ws.Range("M132").Value = -5753.856800000001
$ws.Range("N132").Value = -83527.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18814.094
$ws.Range("I132").Value = 24350.363
$ws.Range("J132").Value = 6634.3
$ws.Range("K132").Value = 73051.08900000001
$ws.Range("L132").Value = 19902.9
$ws.Range("M132").Value = -70521.08900000001
$ws.Range("N132").Value = -24962.9
